# Apply the Sat Oct 12 22:12:45 UTC 2024 "cryptos list" refresh.
# Updates Price (D) / Volume(1h) (E) figures for each coin row, and
# re-sorts the two coin pairs that swapped rank (rows 21/22 and 25/26),
# which also rewrites their Coin (B) and Link (C) cells.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 4).Value = "63.231.74"
$ws.Cells.Item(2, 5).Value = "  +0.69%  "
# Row 3
$ws.Cells.Item(3, 4).Value = "2.480.94"
$ws.Cells.Item(3, 5).Value = "  +2.87%  "
# Row 4
$ws.Cells.Item(4, 5).Value = "  -0.39%  "
# Row 5
$ws.Cells.Item(5, 4).NumberFormat = "@"  # keep Price as text, not a number
$ws.Cells.Item(5, 4).Value = "577.77"
$ws.Cells.Item(5, 5).Value = "  +0.51%  "
# Row 6
$ws.Cells.Item(6, 4).NumberFormat = "@"  # keep Price as text, not a number
$ws.Cells.Item(6, 4).Value = "146.78"
$ws.Cells.Item(6, 5).Value = "  +0.56%  "
# Row 7
$ws.Cells.Item(7, 5).Value = "  +0.19%  "
# Row 8
$ws.Cells.Item(8, 5).Value = "  -0.27%  "
# Row 9
$ws.Cells.Item(9, 4).Value = "2.479.37"
$ws.Cells.Item(9, 5).Value = "  +1.59%  "
# Row 10
$ws.Cells.Item(10, 5).Value = "  +0.27%  "
# Row 11
$ws.Cells.Item(11, 5).Value = "  +1.65%  "
# Row 12
$ws.Cells.Item(12, 4).NumberFormat = "@"  # keep Price as text, not a number
$ws.Cells.Item(12, 4).Value = "5.27"
$ws.Cells.Item(12, 5).Value = "  +0.49%  "
# Row 13
$ws.Cells.Item(13, 5).Value = "  +0.10%  "
# Row 14
$ws.Cells.Item(14, 4).NumberFormat = "@"  # keep Price as text, not a number
$ws.Cells.Item(14, 4).Value = "28.66"
$ws.Cells.Item(14, 5).Value = "  +4.96%  "
# Row 15
$ws.Cells.Item(15, 4).NumberFormat = "@"  # keep Price as text, not a number
$ws.Cells.Item(15, 4).Value = "0.0000179"
$ws.Cells.Item(15, 5).Value = "  +1.17%  "
# Row 16
$ws.Cells.Item(16, 4).Value = "2.931.60"
$ws.Cells.Item(16, 5).Value = "  +2.50%  "
# Row 17
$ws.Cells.Item(17, 4).Value = "63.158.57"
$ws.Cells.Item(17, 5).Value = "  +0.78%  "
# Row 18
$ws.Cells.Item(18, 4).Value = "2.478.89"
$ws.Cells.Item(18, 5).Value = "  +1.45%  "
# Row 19
$ws.Cells.Item(19, 4).NumberFormat = "@"  # keep Price as text, not a number
$ws.Cells.Item(19, 4).Value = "8.27"
$ws.Cells.Item(19, 5).Value = "  +4.35%  "
# Row 20
$ws.Cells.Item(20, 4).NumberFormat = "@"  # keep Price as text, not a number
$ws.Cells.Item(20, 4).Value = "11.05"
$ws.Cells.Item(20, 5).Value = "  +0.73%  "
# Row 21
$ws.Cells.Item(21, 2).Value = "BitcoinCash"
$ws.Cells.Item(21, 3).Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Cells.Item(21, 4).NumberFormat = "@"  # keep Price as text, not a number
$ws.Cells.Item(21, 4).Value = "329.71"
$ws.Cells.Item(21, 5).Value = "  +0.40%  "
# Row 22
$ws.Cells.Item(22, 2).Value = "SuiNetwork"
$ws.Cells.Item(22, 3).Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Cells.Item(22, 4).NumberFormat = "@"  # keep Price as text, not a number
$ws.Cells.Item(22, 4).Value = "2.25"
$ws.Cells.Item(22, 5).Value = "  +9.85%  "
# Row 23
$ws.Cells.Item(23, 4).NumberFormat = "@"  # keep Price as text, not a number
$ws.Cells.Item(23, 4).Value = "4.14"
# Row 24
$ws.Cells.Item(24, 4).NumberFormat = "@"  # keep Price as text, not a number
$ws.Cells.Item(24, 4).Value = "0.999"
$ws.Cells.Item(24, 5).Value = "  +0.19%  "
# Row 25
$ws.Cells.Item(25, 2).Value = "Litecoin"
$ws.Cells.Item(25, 3).Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Cells.Item(25, 4).NumberFormat = "@"  # keep Price as text, not a number
$ws.Cells.Item(25, 4).Value = "66.25"
$ws.Cells.Item(25, 5).Value = "  +1.07%  "
# Row 26
$ws.Cells.Item(26, 2).Value = "Bittensor"
$ws.Cells.Item(26, 3).Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Cells.Item(26, 4).NumberFormat = "@"  # keep Price as text, not a number
$ws.Cells.Item(26, 4).Value = "672.16"
$ws.Cells.Item(26, 5).Value = "  +6.86%  "
# Row 27
$ws.Cells.Item(27, 4).NumberFormat = "@"  # keep Price as text, not a number
$ws.Cells.Item(27, 4).Value = "9.72"
# Row 28
$ws.Cells.Item(28, 4).NumberFormat = "@"  # keep Price as text, not a number
$ws.Cells.Item(28, 4).Value = "0.0000100"
$ws.Cells.Item(28, 5).Value = "  +1.03%  "
# Row 29
$ws.Cells.Item(29, 4).Value = "2.620.20"
$ws.Cells.Item(29, 5).Value = "  +3.57%  "
# Row 30
$ws.Cells.Item(30, 5).Value = "  -9.61%  "
# Row 31
$ws.Cells.Item(31, 4).NumberFormat = "@"  # keep Price as text, not a number
$ws.Cells.Item(31, 4).Value = "1.48"
$ws.Cells.Item(31, 5).Value = "  +3.92%  "
# Row 32
$ws.Cells.Item(32, 4).NumberFormat = "@"  # keep Price as text, not a number
$ws.Cells.Item(32, 4).Value = "8.08"
$ws.Cells.Item(32, 5).Value = "  -1.58%  "
# Row 33
$ws.Cells.Item(33, 5).Value = "  +1.01%  "
# Row 34
$ws.Cells.Item(34, 4).NumberFormat = "@"  # keep Price as text, not a number
$ws.Cells.Item(34, 4).Value = "0.133"
$ws.Cells.Item(34, 5).Value = "  -3.75%  "
# Row 35
$ws.Cells.Item(35, 5).Value = "  +3.80%  "
# Row 36
$ws.Cells.Item(36, 4).NumberFormat = "@"  # keep Price as text, not a number
$ws.Cells.Item(36, 4).Value = "0.999"
$ws.Cells.Item(36, 5).Value = "  +0.34%  "
# Row 37
$ws.Cells.Item(37, 5).Value = "  +0.65%  "
# Row 38
$ws.Cells.Item(38, 5).Value = "  +1.12%  "
# Row 39
$ws.Cells.Item(39, 5).Value = "  -0.63%  "
# Row 40
$ws.Cells.Item(40, 5).Value = "  +0.64%  "
# Row 41
$ws.Cells.Item(41, 4).NumberFormat = "@"  # keep Price as text, not a number
$ws.Cells.Item(41, 4).Value = "151.71"
$ws.Cells.Item(41, 5).Value = "  -0.06%  "
# Row 42
$ws.Cells.Item(42, 4).NumberFormat = "@"  # keep Price as text, not a number
$ws.Cells.Item(42, 4).Value = "2.72"
$ws.Cells.Item(42, 5).Value = "  -2.02%  "
# Row 43
$ws.Cells.Item(43, 4).NumberFormat = "@"  # keep Price as text, not a number
$ws.Cells.Item(43, 4).Value = "1.77"
$ws.Cells.Item(43, 5).Value = "  +0.21%  "
# Row 45
$ws.Cells.Item(45, 4).Value = "0.0₆0313"
$ws.Cells.Item(45, 5).Value = "  -33.43%  "
# Row 46
$ws.Cells.Item(46, 4).NumberFormat = "@"  # keep Price as text, not a number
$ws.Cells.Item(46, 4).Value = "155.66"
$ws.Cells.Item(46, 5).Value = "  +7.11%  "
# Row 47
$ws.Cells.Item(47, 5).Value = "  +13.04%  "
# Row 48
$ws.Cells.Item(48, 5).Value = "  +0.08%  "
# Row 49
$ws.Cells.Item(49, 4).NumberFormat = "@"  # keep Price as text, not a number
$ws.Cells.Item(49, 4).Value = "20.54"
$ws.Cells.Item(49, 5).Value = "  +0.06%  "
# Row 50
$ws.Cells.Item(50, 5).Value = "  +1.22%  "
# Row 51
$ws.Cells.Item(51, 4).NumberFormat = "@"  # keep Price as text, not a number
$ws.Cells.Item(51, 4).Value = "0.0514"
$ws.Cells.Item(51, 5).Value = "  -0.54%  "
